$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("doFindStores")

$ws.Range("A2").Value = "Philadelphia"
$ws.Range("A3").Value = "New york"
$ws.Range("A4").Value = "Washington D.C"
$ws.Range("A5").Value = "Kenwood"

$ws.Range("A7").Select()
